$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.010.07"
$ws.Range("E2").Value = "  -5.50%  "

$ws.Range("D3").Value = "1.820.82"
$ws.Range("E3").Value = "  -5.27%  "

$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  -0.76%  "

$ws.Range("D5").Value = "'328.47"
$ws.Range("E5").Value = "  -3.05%  "

$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  -0.79%  "

$ws.Range("D7").Value = "'0.4611"
$ws.Range("E7").Value = "  -4.16%  "

$ws.Range("D8").Value = "'0.3842"
$ws.Range("E8").Value = "  -5.22%  "

$ws.Range("D9").Value = "'46.03"
$ws.Range("E9").Value = "  -3.85%  "

$ws.Range("D10").Value = "'0.07806"
$ws.Range("E10").Value = "  -3.51%  "

$ws.Range("D11").Value = "'0.9567"
$ws.Range("E11").Value = "  -4.25%  "

$ws.Range("D12").Value = "'21.79"
$ws.Range("E12").Value = "  -6.81%  "

$ws.Range("D13").Value = "'5.657"
$ws.Range("E13").Value = "  -5.41%  "

$ws.Range("D14").Value = "'6.854"
$ws.Range("E14").Value = "  -4.59%  "

$ws.Range("D15").Value = "1.735.20"
$ws.Range("E15").Value = "  -10.77%  "

$ws.Range("D16").Value = "'0.06862"
$ws.Range("E16").Value = "  +0.49%  "

$ws.Range("D17").Value = "'1.005"
$ws.Range("E17").Value = "  -0.78%  "

$ws.Range("D18").Value = "'86.24"
$ws.Range("E18").Value = "  -4.21%  "

$ws.Range("D19").Value = "'0.000009909"
$ws.Range("E19").Value = "  -3.63%  "

$ws.Range("D20").Value = "'16.76"
$ws.Range("E20").Value = "  -4.26%  "

$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  -1.02%  "

$ws.Range("D22").Value = "28.032.35"
$ws.Range("E22").Value = "  -5.47%  "

$ws.Range("D23").Value = "'5.313"
$ws.Range("E23").Value = "  -4.23%  "

$ws.Range("D24").Value = "'10.96"
$ws.Range("E24").Value = "  -6.84%  "

$ws.Range("D25").Value = "'2.124"
$ws.Range("E25").Value = "  -1.63%  "

$ws.Range("D26").Value = "1.985.01"
$ws.Range("E26").Value = "  -9.29%  "

$ws.Range("D27").Value = "'151.71"
$ws.Range("E27").Value = "  -3.54%  "

$ws.Range("D28").Value = "'19.14"
$ws.Range("E28").Value = "  -3.57%  "

$ws.Range("D29").Value = "'5.689"
$ws.Range("E29").Value = "  -13.80%  "

$ws.Range("D30").Value = "'1.969"
$ws.Range("E30").Value = "  -4.60%  "

$ws.Range("D31").Value = "'116.27"
$ws.Range("E31").Value = "  -3.36%  "

$ws.Range("D32").Value = "'0.09267"
$ws.Range("E32").Value = "  -3.35%  "

$ws.Range("D33").Value = "'0.9357"
$ws.Range("E33").Value = "  -6.47%  "

$ws.Range("D34").Value = "'5.255"
$ws.Range("E34").Value = "  -4.86%  "

$ws.Range("D35").Value = "'3.425"
$ws.Range("E35").Value = "  -3.41%  "

$ws.Range("D36").Value = "'1.304"
$ws.Range("E36").Value = "  -6.65%  "

$ws.Range("D37").Value = "'0.05962"
$ws.Range("E37").Value = "  -8.87%  "

$ws.Range("E38").Value = "  -5.19%  "

$ws.Range("D39").Value = "'1.147"
$ws.Range("E39").Value = "  -4.34%  "

$ws.Range("E40").Value = "  -0.99%  "

$ws.Range("D41").Value = "'7.531"
$ws.Range("E41").Value = "  -4.11%  "

$ws.Range("D42").Value = "'0.5584"
$ws.Range("E42").Value = "  -5.13%  "

$ws.Range("D43").Value = "'9.891"
$ws.Range("E43").Value = "  -7.04%  "

$ws.Range("E44").Value = "  -3.54%  "

$ws.Range("E45").Value = "  -1.68%  "

$ws.Range("D46").Value = "'2.224"
$ws.Range("E46").Value = "  -10.31%  "

$ws.Range("D47").Value = "'11.57"
$ws.Range("E47").Value = "  -5.79%  "

$ws.Range("D48").Value = "'0.5239"
$ws.Range("E48").Value = "  -4.86%  "

$ws.Range("D49").Value = "'0.06987"
$ws.Range("E49").Value = "  -6.38%  "

$ws.Range("D50").Value = "'1.823"
$ws.Range("E50").Value = "  -6.98%  "

$ws.Range("E51").Value = "  -3.52%  "
